$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 728.7692
$ws.Cells.Item(92, 9).Value = 553.2222
$ws.Cells.Item(92, 11).Value = 553.2222
$ws.Cells.Item(92, 13).Value = 694.7778
$ws.Cells.Item(98, 8).Value = 859.1
$ws.Cells.Item(98, 9).Value = 1057.1578
$ws.Cells.Item(98, 10).Value = 517
$ws.Cells.Item(98, 11).Value = 1057.1578
$ws.Cells.Item(98, 12).Value = 517
$ws.Cells.Item(98, 13).Value = 440.8422
$ws.Cells.Item(98, 14).Value = -3513
$ws.Cells.Item(107, 8).Value = 683.63635
$ws.Cells.Item(107, 9).Value = 706.1667
$ws.Cells.Item(107, 11).Value = 706.1667
$ws.Cells.Item(107, 13).Value = 1213.8333
$ws.Cells.Item(111, 8).Value = 991.6667
$ws.Cells.Item(111, 10).Value = 887
$ws.Cells.Item(111, 12).Value = 2661
$ws.Cells.Item(111, 14).Value = -8795
$ws.Cells.Item(122, 8).Value = 859.1
$ws.Cells.Item(122, 9).Value = 1057.1578
$ws.Cells.Item(122, 10).Value = 517
$ws.Cells.Item(122, 11).Value = 3171.4734
$ws.Cells.Item(122, 12).Value = 1551
$ws.Cells.Item(122, 13).Value = -721.4733999999999
$ws.Cells.Item(122, 14).Value = -6451
$ws.Cells.Item(138, 8).Value = 4618.7
$ws.Cells.Item(138, 10).Value = 4946.1577
$ws.Cells.Item(138, 12).Value = 14838.4731
$ws.Cells.Item(138, 14).Value = -25118.4731
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(140, 14).ClearContents()
$ws.Cells.Item(141, 8).Value = 8999.5
$ws.Cells.Item(141, 9).Value = 8000
$ws.Cells.Item(141, 11).Value = 24000
$ws.Cells.Item(141, 13).Value = -18820
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2992.3572
$ws.Cells.Item(45, 9).Value = 2598
$ws.Cells.Item(45, 11).Value = 2598
$ws.Cells.Item(45, 13).Value = -2221
$ws.Cells.Item(74, 8).Value = 1439.6154
$ws.Cells.Item(74, 9).Value = 966
$ws.Cells.Item(74, 11).Value = 966
$ws.Cells.Item(74, 13).Value = -92
$ws.Cells.Item(77, 8).Value = 1439.6154
$ws.Cells.Item(77, 9).Value = 966
$ws.Cells.Item(77, 11).Value = 4830
$ws.Cells.Item(77, 13).Value = -462
$ws.Cells.Item(110, 8).Value = 1149.1666
$ws.Cells.Item(110, 9).Value = 1149.1666
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 11).Value = 1149.1666
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 13).Value = 895.8334
$ws.Cells.Item(110, 14).ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 1330
$ws.Cells.Item(64, 10).Value = 1495
$ws.Cells.Item(64, 12).Value = 1495
$ws.Cells.Item(64, 14).Value = -1945
$ws.Cells.Item(67, 8).Value = 1330
$ws.Cells.Item(67, 10).Value = 1495
$ws.Cells.Item(67, 12).Value = 1495
$ws.Cells.Item(67, 14).Value = -3055
$ws.Cells.Item(134, 8).Value = 2857.16
$ws.Cells.Item(134, 9).Value = 2654.4736
$ws.Cells.Item(134, 11).Value = 7963.4208
$ws.Cells.Item(134, 13).Value = -5428.4208
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2153.1304
$ws.Cells.Item(31, 9).Value = 2177.75
$ws.Cells.Item(31, 11).Value = 2177.75
$ws.Cells.Item(31, 13).Value = -1882.75
$ws.Cells.Item(34, 8).Value = 2153.1304
$ws.Cells.Item(34, 9).Value = 2177.75
$ws.Cells.Item(34, 11).Value = 2177.75
$ws.Cells.Item(34, 13).Value = -1975.75
$ws.Cells.Item(69, 8).Value = 16500
$ws.Cells.Item(69, 9).Value = 5000
$ws.Cells.Item(69, 10).Value = 20333.334
$ws.Cells.Item(69, 11).Value = 5000
$ws.Cells.Item(69, 12).Value = 20333.334
$ws.Cells.Item(69, 13).Value = -4251
$ws.Cells.Item(69, 14).Value = -21831.334
$ws.Cells.Item(72, 8).Value = 16500
$ws.Cells.Item(72, 9).Value = 5000
$ws.Cells.Item(72, 10).Value = 20333.334
$ws.Cells.Item(72, 11).Value = 15000
$ws.Cells.Item(72, 12).Value = 61000.00199999999
$ws.Cells.Item(72, 13).Value = -11256
$ws.Cells.Item(72, 14).Value = -68488.00199999999
$ws.Cells.Item(105, 8).Value = 3655.2
$ws.Cells.Item(105, 9).Value = 3319
$ws.Cells.Item(105, 11).Value = 3319
$ws.Cells.Item(105, 13).Value = -1572
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 270.25
$ws.Cells.Item(7, 10).Value = 300
$ws.Cells.Item(7, 12).Value = 900
$ws.Cells.Item(7, 14).Value = -1124
$ws.Cells.Item(75, 8).Value = 2010.6666
$ws.Cells.Item(75, 10).Value = 2010.6666
$ws.Cells.Item(75, 12).Value = 6031.9998
$ws.Cells.Item(75, 14).Value = -8027.9998
$ws.Cells.Item(78, 8).Value = 2010.6666
$ws.Cells.Item(78, 10).Value = 2010.6666
$ws.Cells.Item(78, 12).Value = 18095.9994
$ws.Cells.Item(78, 14).Value = -28079.9994
$ws.Cells.Item(92, 8).Value = 1663.3334
$ws.Cells.Item(92, 9).Value = 1663.3334
$ws.Cells.Item(92, 11).Value = 4990.0002
$ws.Cells.Item(92, 13).Value = -3742.0002
$ws.Cells.Item(94, 8).Value = 2779.8
$ws.Cells.Item(94, 10).Value = 2975
$ws.Cells.Item(94, 12).Value = 8925
$ws.Cells.Item(94, 14).Value = -10277
$ws.Cells.Item(138, 8).Value = 4674.75
$ws.Cells.Item(138, 9).Value = 4674.75
$ws.Cells.Item(138, 11).Value = 14024.25
$ws.Cells.Item(138, 13).Value = -8884.25
$ws.Cells.Item(139, 8).Value = 1397.375
$ws.Cells.Item(139, 9).Value = 1397.375
$ws.Cells.Item(139, 11).Value = 4192.125
$ws.Cells.Item(139, 13).Value = 947.875
$ws.Cells.Item(140, 8).Value = 1618.1428
$ws.Cells.Item(140, 9).Value = 1618.1428
$ws.Cells.Item(140, 11).Value = 4854.428400000001
$ws.Cells.Item(140, 13).Value = 325.5715999999993
$ws.Cells.Item(141, 8).Value = 5740.2856
$ws.Cells.Item(141, 9).Value = 5740.2856
$ws.Cells.Item(141, 11).Value = 17220.8568
$ws.Cells.Item(141, 13).Value = -12040.8568
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 65.71429000000001
$ws.Cells.Item(2, 9).Value = 72.5
$ws.Cells.Item(2, 11).Value = 72.5
$ws.Cells.Item(2, 13).Value = 40.5
$ws.Cells.Item(46, 8).Value = 4324.2856
$ws.Cells.Item(46, 10).Value = 4324.2856
$ws.Cells.Item(46, 12).Value = 4324.2856
$ws.Cells.Item(46, 14).Value = -4636.2856
$ws.Cells.Item(97, 8).Value = 793.875
$ws.Cells.Item(97, 9).Value = 665.4211
$ws.Cells.Item(97, 10).Value = 1282
$ws.Cells.Item(97, 11).Value = 665.4211
$ws.Cells.Item(97, 12).Value = 1282
$ws.Cells.Item(97, 13).Value = -169.4211
$ws.Cells.Item(97, 14).Value = -2274
$ws.Cells.Item(122, 8).Value = 4659.8
$ws.Cells.Item(122, 9).Value = 3824.75
$ws.Cells.Item(122, 11).Value = 11474.25
$ws.Cells.Item(122, 13).Value = -9024.25
$ws.Cells.Item(126, 8).Value = 4056.5
$ws.Cells.Item(126, 9).Value = 3410.25
$ws.Cells.Item(126, 11).Value = 10230.75
$ws.Cells.Item(126, 13).Value = -7760.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 7616.6665
$ws.Cells.Item(136, 8).Value = 5592.2354
$ws.Cells.Item(136, 9).Value = 5647.9287
$ws.Cells.Item(136, 11).Value = 16943.7861
$ws.Cells.Item(136, 13).Value = -14393.7861
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 7577.3335
$ws.Cells.Item(81, 9).Value = 4269.6665
$ws.Cells.Item(81, 10).Value = 17500.334
$ws.Cells.Item(81, 11).Value = 8539.333000000001
$ws.Cells.Item(81, 12).Value = 35000.668
$ws.Cells.Item(81, 13).Value = -7478.333000000001
$ws.Cells.Item(81, 14).Value = -37122.668
$ws.Cells.Item(84, 8).Value = 7577.3335
$ws.Cells.Item(84, 9).Value = 4269.6665
$ws.Cells.Item(84, 10).Value = 17500.334
$ws.Cells.Item(84, 11).Value = 42696.665
$ws.Cells.Item(84, 12).Value = 175003.34
$ws.Cells.Item(84, 13).Value = -37392.665
$ws.Cells.Item(84, 14).Value = -185611.34
$ws.Cells.Item(107, 8).Value = 557.7273
$ws.Cells.Item(107, 9).Value = 512.375
$ws.Cells.Item(107, 11).Value = 1537.125
$ws.Cells.Item(107, 13).Value = 382.875
$ws.Cells.Item(132, 8).Value = 1911.4445
$ws.Cells.Item(132, 9).Value = 1393.9166
$ws.Cells.Item(132, 11).Value = 4181.7498
$ws.Cells.Item(132, 13).Value = -1651.7498
